{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the target paragraphs by their (current) text content so the\n// script is resilient to exact index assumptions.\nlet emptyListPara = null;   // empty List Paragraph item right before \"Invite...\"\nlet piazzaPara = null;      // \"Invite everyone to Piazza...\" paragraph\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text === \"\" && i + 1 < paragraphs.items.length &&\n      paragraphs.items[i + 1].text.indexOf(\"Invite everyone to Piazza\") !== -1) {\n    emptyListPara = para;\n  }\n  if (para.text.indexOf(\"Invite everyone to Piazza\") !== -1) {\n    piazzaPara = para;\n  }\n}\n\nif (!emptyListPara || !piazzaPara) {\n  throw new Error(\"Could not locate the expected paragraphs to edit.\");\n}\n\n// 1) Turn the previously-empty, list-numbered paragraph into a plain\n//    \"Normal\" paragraph (no list, no ListParagraph style) carrying the\n//    new sentence about Meeting Times, but keep the 6pt (120 twip)\n//    space-before.\nemptyListPara.style = \"Normal\";\nemptyListPara.spaceBefore = 6;\nemptyListPara.insertText(\"Have them fill out the Meeting Times (TENTATIVE).\", Word.InsertLocation.start);\n\n// 2) Give the existing \"Invite everyone to Piazza...\" paragraph the same\n//    6pt space-before (it previously had no paragraph properties at all).\npiazzaPara.spaceBefore = 6;\n\n// 3) Insert a brand-new paragraph right after it for the new closing note.\nconst newPara = piazzaPara.insertParagraph(\n  \"Let them do Computer Setup and, perhaps, some of all of Unit 1.\",\n  Word.InsertLocation.after\n);\nnewPara.spaceBefore = 6;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Invite everyone to Piazza...\" paragraph by its text so the\n# script does not depend on a brittle, hard-coded paragraph index.\n$piazzaIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*Invite everyone to Piazza*\") {\n        $piazzaIdx = $i\n    }\n}\n\nif ($piazzaIdx -eq -1) {\n    throw \"Could not locate the 'Invite everyone to Piazza' paragraph.\"\n}\n\n$piazzaPara = $d.Paragraphs.Item($piazzaIdx)\n$emptyListPara = $d.Paragraphs.Item($piazzaIdx - 1)\n\n# 1) Turn the previously-empty, list-numbered paragraph into a plain\n#    \"Normal\" paragraph (no list, no ListParagraph style) carrying the\n#    new sentence about Meeting Times, but keep the 6pt (120 twip)\n#    space-before.\n$emptyListPara.Style = \"Normal\"\n$emptyListPara.SpaceBefore = 6\n$emptyListPara.Range.Text = \"Have them fill out the Meeting Times (TENTATIVE).\"\n\n# 2) Give the existing \"Invite everyone to Piazza...\" paragraph the same\n#    6pt space-before (it previously had no paragraph properties at all).\n$piazzaPara.SpaceBefore = 6\n\n# 3) Insert a brand-new paragraph right after it for the new closing note.\n$piazzaPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs.Item($piazzaIdx + 1)\n$newPara.Range.Text = \"Let them do Computer Setup and, perhaps, some of all of Unit 1.\"\n$newPara.SpaceBefore = 6\n"}
